$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 (I0) and J1 (IF), matching style of existing header H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate I and J columns with numeric data for rows 2-62
$iValues = @(8,7,2,3,9,5,6,7,2,7,8,6,8,8,7,7,9,7,6,7,8,6,8,7,8,5,5,8,7,6,7,10,7,2,8,5,5,7,7,7,6,5,6,7,7,9,9,7,9,8,4,6,7,9,10,6,6,4,8,5,4)
$jValues = @(8,7,3,5,9,7,6,8,3,7,8,8,8,8,7,7,9,8,8,7,9,6,8,7,8,5,6,8,7,7,7,10,7,2,8,6,6,8,8,8,6,5,6,7,7,9,9,8,9,8,5,6,7,9,10,6,6,6,9,5,4)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $r = $idx + 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
